# "Basic framework for powerpoint added"
#
# - Slide 4 ("Title and Content") is rebuilt as a "Picture with Caption"
#   slide: title text is set, the old Content Placeholder is swapped for a
#   Picture Placeholder + a half-width Text Placeholder.
# - Two new "Title and Content" slides are appended (Class Explanation,
#   Feedback Review).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 4: "Method/Function Explanation" (Picture with Caption layout)
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)

# Grab a throw-away slide using the "Picture with Caption" layout (36) so
# we can hand its CustomLayout to the real slide 4 and then discard it.
$scratch = $p.Slides.Add($p.Slides.Count + 1, 36)
$pictureLayout = $scratch.CustomLayout

# Drop the old body placeholder - it has no equivalent on the new layout -
# then re-point slide 4 at the picture-with-caption layout, which adds the
# Picture Placeholder + half Text Placeholder shapes.
$s4.Shapes.Item(2).Delete()
$s4.CustomLayout = $pictureLayout

$scratch.Delete()

$title4 = $s4.Shapes.Item(1)
$title4.TextFrame.TextRange.Text = "Method/Function Explanation"
$title4.TextFrame.TextRange.LanguageID = "en-AU"

$picPh = $s4.Shapes.Item(2)
$picPh.Name = "Picture Placeholder 3"

$textPh = $s4.Shapes.Item(3)
$textPh.Name = "Text Placeholder 4"
$textPh.TextFrame.TextRange.Text = "Picking function that has imported code in it"
$textPh.TextFrame.TextRange.LanguageID = "en-AU"

# ---------------------------------------------------------------------
# Slide 5: "Class Explanation" (Title and Content layout)
# ---------------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Class Explanation"
$s5.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-AU"

# ---------------------------------------------------------------------
# Slide 6: "Feedback Review" (Title and Content layout)
# ---------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)
$s6.Shapes.Item(1).TextFrame.TextRange.Text = "Feedback Review"
$s6.Shapes.Item(1).TextFrame.TextRange.LanguageID = "en-AU"
